$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Ty le %" (contribution percentage) column for each team
# member with the grades the team assigned to themselves.
$ws.Range("C3").Value = 0.9
$ws.Range("C4").Value = 0.89
$ws.Range("C5").Value = 0.97
$ws.Range("C6").Value = 0.97
$ws.Range("C7").Value = 0.99

# Display the new values as percentages.
$ws.Range("C3:C7").NumberFormat = "0%"

# Leave the cursor where the user finished editing.
$ws.Range("F7").Select()
